# Update Leve profit-tracking figures (currentAveragePrice / LevePrice / LeveProfit
# columns H:N) for several rows across the ALC, CRP, CUL, GSM, LTW and WVR sheets,
# reflecting refreshed market-board pricing data from the scheduled runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 42662.68
$ws.Range("I64").Value = 145238.58
$ws.Range("J64").Value = 2772.0557
$ws.Range("K64").Value = 145238.58
$ws.Range("L64").Value = 2772.0557
$ws.Range("M64").Value = -144990.58
$ws.Range("N64").Value = -3268.0557
$ws.Range("H67").Value = 42662.68
$ws.Range("I67").Value = 145238.58
$ws.Range("J67").Value = 2772.0557
$ws.Range("K67").Value = 145238.58
$ws.Range("L67").Value = 2772.0557
$ws.Range("M67").Value = -144380.58
$ws.Range("N67").Value = -4488.0557
$ws.Range("H86").Value = 16708.908
$ws.Range("I86").Value = 15286.75
$ws.Range("J86").Value = 20501.334
$ws.Range("K86").Value = 15286.75
$ws.Range("L86").Value = 20501.334
$ws.Range("M86").Value = -14163.75
$ws.Range("N86").Value = -22747.334
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").Value = ""
$ws.Range("H89").Value = 16708.908
$ws.Range("I89").Value = 15286.75
$ws.Range("J89").Value = 20501.334
$ws.Range("K89").Value = 76433.75
$ws.Range("L89").Value = 102506.67
$ws.Range("M89").Value = -70817.75
$ws.Range("N89").Value = -113738.67
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").Value = ""
$ws.Range("H107").Value = 907.5455
$ws.Range("I107").Value = 606.1539
$ws.Range("J107").Value = 1342.8889
$ws.Range("K107").Value = 606.1539
$ws.Range("L107").Value = 1342.8889
$ws.Range("M107").Value = 1313.8461
$ws.Range("N107").Value = -5182.8889
$ws.Range("H133").Value = 73435.46
$ws.Range("J133").Value = 73435.46
$ws.Range("L133").Value = 73435.46
$ws.Range("N133").Value = -83555.46
$ws.Range("H134").Value = 38650
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 38650
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 38650
$ws.Range("M134").Value = ""
$ws.Range("N134").Value = -48790
$ws.Range("H137").Value = 2561.2654
$ws.Range("I137").Value = 831.8485
$ws.Range("J137").Value = 3439.2769
$ws.Range("K137").Value = 2495.5455
$ws.Range("L137").Value = 10317.8307
$ws.Range("M137").Value = 54.45450000000028
$ws.Range("N137").Value = -15417.8307
$ws.Range("H141").Value = 2922.3704
$ws.Range("I141").Value = 1361.381
$ws.Range("J141").Value = 8385.833
$ws.Range("K141").Value = 4084.143
$ws.Range("L141").Value = 25157.499
$ws.Range("M141").Value = 1095.857
$ws.Range("N141").Value = -35517.499

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 12525
$ws.Range("I11").Value = 50
$ws.Range("J11").Value = 25000
$ws.Range("K11").Value = 50
$ws.Range("L11").Value = 25000
$ws.Range("M11").Value = 90
$ws.Range("N11").Value = -25280
$ws.Range("H31").Value = 4378.9326
$ws.Range("I31").Value = 1691.5946
$ws.Range("J31").Value = 7066.2705
$ws.Range("K31").Value = 1691.5946
$ws.Range("L31").Value = 7066.2705
$ws.Range("M31").Value = -1396.5946
$ws.Range("N31").Value = -7656.2705
$ws.Range("H34").Value = 4378.9326
$ws.Range("I34").Value = 1691.5946
$ws.Range("J34").Value = 7066.2705
$ws.Range("K34").Value = 1691.5946
$ws.Range("L34").Value = 7066.2705
$ws.Range("M34").Value = -1489.5946
$ws.Range("N34").Value = -7470.2705
$ws.Range("H58").Value = 2171.9075
$ws.Range("I58").Value = 1857.2128
$ws.Range("J58").Value = 4284.857
$ws.Range("K58").Value = 1857.2128
$ws.Range("L58").Value = 4284.857
$ws.Range("M58").Value = -1654.2128
$ws.Range("N58").Value = -4690.857
$ws.Range("I62").Value = 2716.6667
$ws.Range("J62").Value = 2916.6667
$ws.Range("K62").Value = 2716.6667
$ws.Range("L62").Value = 2916.6667
$ws.Range("M62").Value = -2092.6667
$ws.Range("N62").Value = -4164.6667
$ws.Range("I65").Value = 2716.6667
$ws.Range("J65").Value = 2916.6667
$ws.Range("K65").Value = 13583.3335
$ws.Range("L65").Value = 14583.3335
$ws.Range("M65").Value = -10463.3335
$ws.Range("N65").Value = -20823.3335
$ws.Range("H132").Value = 34683.44
$ws.Range("I132").Value = 1685.931
$ws.Range("J132").Value = 103035.43
$ws.Range("K132").Value = 5057.793
$ws.Range("L132").Value = 309106.29
$ws.Range("M132").Value = -2527.793
$ws.Range("N132").Value = -314166.29
$ws.Range("H134").Value = 2216.182
$ws.Range("I134").Value = 1128.125
$ws.Range("J134").Value = 5117.6665
$ws.Range("K134").Value = 3384.375
$ws.Range("L134").Value = 15352.9995
$ws.Range("M134").Value = -849.375
$ws.Range("N134").Value = -20422.9995
$ws.Range("H136").Value = 2171.9075
$ws.Range("I136").Value = 1857.2128
$ws.Range("J136").Value = 4284.857
$ws.Range("K136").Value = 5571.6384
$ws.Range("L136").Value = 12854.571
$ws.Range("M136").Value = -3021.6384
$ws.Range("N136").Value = -17954.571
$ws.Range("H141").Value = 19445
$ws.Range("J141").Value = 19890
$ws.Range("L141").Value = 19890
$ws.Range("N141").Value = -30250

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 2706.6667
$ws.Range("I11").Value = 2975
$ws.Range("K11").Value = 8925
$ws.Range("M11").Value = -8785
$ws.Range("H131").Value = 902.4141
$ws.Range("I131").Value = 658
$ws.Range("J131").Value = 918.1828
$ws.Range("K131").Value = 1974
$ws.Range("L131").Value = 2754.5484
$ws.Range("M131").Value = 3066
$ws.Range("N131").Value = -12834.5484

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 361642.84
$ws.Range("I80").Value = 630100
$ws.Range("J80").Value = 3700
$ws.Range("K80").Value = 630100
$ws.Range("L80").Value = 3700
$ws.Range("M80").Value = -629102
$ws.Range("N80").Value = -5696
$ws.Range("H83").Value = 361642.84
$ws.Range("I83").Value = 630100
$ws.Range("J83").Value = 3700
$ws.Range("K83").Value = 3150500
$ws.Range("L83").Value = 18500
$ws.Range("M83").Value = -3145508
$ws.Range("N83").Value = -28484
$ws.Range("H132").Value = 2066.913
$ws.Range("I132").Value = 1643.5312
$ws.Range("J132").Value = 3034.6428
$ws.Range("K132").Value = 4930.5936
$ws.Range("L132").Value = 9103.9284
$ws.Range("M132").Value = -2400.5936
$ws.Range("N132").Value = -14163.9284

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1998
$ws.Range("I136").Value = 1699.1111
$ws.Range("J136").Value = 3612
$ws.Range("K136").Value = 5097.3333
$ws.Range("L136").Value = 10836
$ws.Range("M136").Value = -2547.3333
$ws.Range("N136").Value = -15936

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1571.6875
$ws.Range("I132").Value = 1141.6666
$ws.Range("J132").Value = 2861.75
$ws.Range("K132").Value = 3424.9998
$ws.Range("L132").Value = 8585.25
$ws.Range("M132").Value = -894.9998
$ws.Range("N132").Value = -13645.25

Write-Host "Masamune_Profits: refreshed pricing figures applied."
